# Update the three-digit-by-one-digit division problems in the table.
# Each old expression is unique in the document, so a simple
# Find/Replace (not "Replace All" needed, but harmless either way)
# for each pair safely updates exactly the intended cell.

$d = $word.ActiveDocument

$pairs = @(
    @("460÷2=", "623÷9="),
    @("642÷6=", "522÷7="),
    @("788÷3=", "549÷7="),
    @("114÷9=", "465÷7="),
    @("432÷3=", "451÷2="),
    @("307÷8=", "237÷2="),
    @("493÷3=", "731÷8="),
    @("469÷3=", "485÷5="),
    @("696÷5=", "259÷5="),
    @("176÷4=", "277÷9="),
    @("388÷9=", "327÷9="),
    @("467÷7=", "765÷6="),
    @("799÷8=", "128÷2="),
    @("836÷3=", "347÷9="),
    @("351÷2=", "571÷7="),
    @("224÷5=", "627÷4="),
    @("298÷6=", "337÷9="),
    @("769÷4=", "928÷4="),
    @("293÷9=", "231÷6="),
    @("116÷2=", "587÷8="),
    @("642÷8=", "549÷9="),
    @("194÷3=", "756÷5="),
    @("227÷6=", "612÷5="),
    @("973÷4=", "910÷6="),
    @("454÷5=", "447÷8=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
